$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1654135338345865
$ws.Range("C2").Value = 0.6165413533834586
$ws.Range("J2").Value = 0.03007518796992481
$ws.Range("P2").Value = 0.1203007518796992
$ws.Range("S2").Value = 0.06766917293233082
$ws.Range("B3").Value = 0.01219512195121951
$ws.Range("C3").Value = 0.02439024390243903
$ws.Range("J3").Value = 0.02439024390243903
$ws.Range("P3").Value = 0.7317073170731707
$ws.Range("S3").Value = 0.2073170731707317
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.03508771929824561
$ws.Range("D6").Value = 0.02631578947368421
$ws.Range("F6").Value = 0.1140350877192982
$ws.Range("J6").Value = 0.2807017543859649
$ws.Range("O6").Value = 0.008771929824561403
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.03508771929824561
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1066666666666667
$ws.Range("F7").Value = 0.05333333333333334
$ws.Range("J7").Value = 0.1733333333333333
$ws.Range("O7").Value = 0.02666666666666667
$ws.Range("Q7").Value = 0.2533333333333334
$ws.Range("R7").Value = 0.06666666666666667
$ws.Range("S7").Value = 0.32
$ws.Range("B8").Value = 0.08633093525179857
$ws.Range("D8").Value = 0.02877697841726619
$ws.Range("F8").Value = 0.07194244604316546
$ws.Range("J8").Value = 0.158273381294964
$ws.Range("O8").Value = 0.02877697841726619
$ws.Range("Q8").Value = 0.2446043165467626
$ws.Range("R8").Value = 0.04316546762589928
$ws.Range("S8").Value = 0.3381294964028777
$ws.Range("B9").Value = 0.1168831168831169
$ws.Range("F9").Value = 0.05194805194805195
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("O9").Value = 0.05194805194805195
$ws.Range("Q9").Value = 0.1688311688311688
$ws.Range("R9").Value = 0.07792207792207792
$ws.Range("S9").Value = 0.4415584415584415
$ws.Range("B10").Value = 0.1441605839416058
$ws.Range("D10").Value = 0.02372262773722628
$ws.Range("E10").Value = 0.001824817518248175
$ws.Range("F10").Value = 0.08029197080291971
$ws.Range("J10").Value = 0.08941605839416059
$ws.Range("O10").Value = 0.02554744525547445
$ws.Range("Q10").Value = 0.1934306569343066
$ws.Range("R10").Value = 0.06569343065693431
$ws.Range("S10").Value = 0.3759124087591241
$ws.Range("G11").Value = 0.09565217391304348
$ws.Range("J11").Value = 0.1478260869565217
$ws.Range("K11").Value = 0.1565217391304348
$ws.Range("L11").Value = 0.591304347826087
$ws.Range("S11").Value = 0.008695652173913044
$ws.Range("G12").Value = 0.676056338028169
$ws.Range("J12").Value = 0.2535211267605634
$ws.Range("K12").Value = 0.01408450704225352
$ws.Range("L12").Value = 0.04225352112676056
$ws.Range("S12").Value = 0.01408450704225352
$ws.Range("G13").Value = 0.68
$ws.Range("J13").Value = 0.32
$ws.Range("F15").Value = 0.01886792452830189
$ws.Range("H15").Value = 0.1037735849056604
$ws.Range("I15").Value = 0.05660377358490566
$ws.Range("J15").Value = 0.3679245283018868
$ws.Range("K15").Value = 0.07547169811320754
$ws.Range("O15").Value = 0.1320754716981132
$ws.Range("S15").Value = 0.2452830188679245
$ws.Range("F16").Value = 0.02325581395348837
$ws.Range("H16").Value = 0.1744186046511628
$ws.Range("I16").Value = 0.1162790697674419
$ws.Range("J16").Value = 0.4302325581395349
$ws.Range("K16").Value = 0.03488372093023256
$ws.Range("M16").Value = 0.02325581395348837
$ws.Range("O16").Value = 0.05813953488372093
$ws.Range("S16").Value = 0.1395348837209302
$ws.Range("F17").Value = 0.03684210526315789
$ws.Range("H17").Value = 0.1578947368421053
$ws.Range("I17").Value = 0.06842105263157895
$ws.Range("J17").Value = 0.4052631578947368
$ws.Range("K17").Value = 0.1368421052631579
$ws.Range("M17").Value = 0.005263157894736842
$ws.Range("O17").Value = 0.06842105263157895
$ws.Range("S17").Value = 0.1210526315789474
$ws.Range("F18").Value = 0.03448275862068965
$ws.Range("H18").Value = 0.1206896551724138
$ws.Range("I18").Value = 0.103448275862069
$ws.Range("J18").Value = 0.4827586206896552
$ws.Range("K18").Value = 0.08620689655172414
$ws.Range("M18").Value = 0.05172413793103448
$ws.Range("O18").Value = 0.03448275862068965
$ws.Range("S18").Value = 0.08620689655172414
$ws.Range("F19").Value = 0.0298804780876494
$ws.Range("H19").Value = 0.149402390438247
$ws.Range("I19").Value = 0.08366533864541832
$ws.Range("J19").Value = 0.3964143426294821
$ws.Range("K19").Value = 0.1055776892430279
$ws.Range("M19").Value = 0.03784860557768924
$ws.Range("N19").Value = 0.00398406374501992
$ws.Range("O19").Value = 0.07171314741035857
$ws.Range("S19").Value = 0.1215139442231076
